# Daily update at 8 AM UTC
# Appends the next day's row (row 50) to the "Wins Over Time" tracking sheet:
#   A50 = 45999 (date, same number format as the rows above it)
#   B50 = 116
#   C50 = 125
#   D50 = 116

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = 45999
$ws.Range("B50").Value = 116
$ws.Range("C50").Value = 125
$ws.Range("D50").Value = 116

# Match the date formatting used by the other cells in column A.
$ws.Range("A50").NumberFormat = $ws.Range("A49").NumberFormat
